# Add "Feuer" / "feuer Y" rows (mirroring the Schlange X/Y fire-movement rows)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A23").Value = "Feuer"
$ws.Range("B23").Value = 400
$ws.Range("C23").Value = 100
$ws.Range("D23").Value = 400

$ws.Range("A24").Value = "feuer Y"
$ws.Range("B24").Value = 400
$ws.Range("C24").Value = 150
$ws.Range("D24").Value = 400

$ws.Range("D25").Select()
